$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per the commit diff.
# NumberFormat is forced to text ("@") before assigning values so that
# Excel does not reinterpret numeric-looking strings (e.g. "9.00", "0.0785")
# as numbers and strip formatting / introduce floating point artifacts.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.527.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.479.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.62"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.23"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.545"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.05%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.70"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0785"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.861.94"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.85"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.19"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +9.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.468.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.765"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.527.66"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.29"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.09"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.68%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.87"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.66"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.99"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.02"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.47"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.43%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.89%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.34"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.42%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ApeXProtocol"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -8.06%  "

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.105"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.03%  "

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.91"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.95%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.39%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.15"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.51"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.59%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.96"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.719.31"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.75%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.32%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.90%  "
